$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.034.71"
$ws.Range("E2").Value = "  -4.89%  "

$ws.Range("D3").Value = "3.077.13"
$ws.Range("E3").Value = "  -4.91%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.24%  "

$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").Value = "3.076.23"
$ws.Range("E8").Value = "  -4.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.19%  "

$ws.Range("E10").Value = "  -6.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -11.64%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.464"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.31%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000217"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.13%  "

$ws.Range("D15").Value = "3.569.83"
$ws.Range("E15").Value = "  -5.11%  "

$ws.Range("D16").Value = "62.992.20"
$ws.Range("E16").Value = "  -5.00%  "

$ws.Range("E17").Value = "  -3.17%  "

$ws.Range("D18").Value = "3.068.99"
$ws.Range("E18").Value = "  -5.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "488.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -10.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.708"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.40%  "

$ws.Range("E23").Value = "  -7.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.24%  "

$ws.Range("E28").Value = "  -6.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -13.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "57.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "519.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.85%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -11.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0403"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -14.22%  "

$ws.Range("D39").Value = "3.100.74"
$ws.Range("E39").Value = "  -1.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0798"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.73%  "

$ws.Range("E41").Value = "  -5.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -12.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.255"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.18%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.108"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.98%  "

$ws.Range("D50").Value = "0.0₃0499"
$ws.Range("E50").Value = "  -12.19%  "

$ws.Range("E51").Value = "  +50.15%  "
